$p = $ppt.ActivePresentation

# --- 1. Footer "date" placeholder text: 1/22/2014 -> 1/28/2014 ---
# This cached auto-date field lives on the Slide Master and on every
# Slide Layout (the "Date Placeholder" shape inherited by slides).
$m = $p.SlideMaster

for ($i = 1; $i -le $m.Shapes.Count; $i++) {
    $sh = $m.Shapes.Item($i)
    if ($sh.Name -like "Date Placeholder*") {
        $sh.TextFrame.TextRange.Text = "1/28/2014"
    }
}

for ($i = 1; $i -le $m.CustomLayouts.Count; $i++) {
    $cl = $m.CustomLayouts.Item($i)
    for ($j = 1; $j -le $cl.Shapes.Count; $j++) {
        $sh = $cl.Shapes.Item($j)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = "1/28/2014"
        }
    }
}

# --- 2. "Multiply 13" shape fill color: FF0000 -> C00000 ---
$s = $p.Slides.Item(1)
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Name -eq "Multiply 13") {
        $sh.Fill.ForeColor.RGB = 192
    }
}
